# Commit: "Fruta / hortaliza, semanal"
# Weekly refresh of the Hortaliza (Arveja Verde) price sheet for
# Comercializadora del Agro de Limari: updates Fecha, Volumen,
# Precio minimo/maximo/promedio ponderado and Precio $/Kg per record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44462
$ws.Range("J2").Value = 400
$ws.Range("K2").Value = 22000
$ws.Range("L2").Value = 23000
$ws.Range("M2").Value = 22500
$ws.Range("P2").Value = 900

$ws.Range("D3").Value = 44349
$ws.Range("J3").Value = 600
$ws.Range("K3").Value = 26000
$ws.Range("L3").Value = 28000
$ws.Range("M3").Value = 27000
$ws.Range("P3").Value = 1080

$ws.Range("D4").Value = 44406
$ws.Range("J4").Value = 600
$ws.Range("K4").Value = 26000
$ws.Range("L4").Value = 28000
$ws.Range("M4").Value = 27000
$ws.Range("P4").Value = 1080

$ws.Range("D5").Value = 44419
$ws.Range("J5").Value = 600
$ws.Range("K5").Value = 27000
$ws.Range("L5").Value = 29000
$ws.Range("M5").Value = 28000
$ws.Range("P5").Value = 1120

$ws.Range("D6").Value = 44475
$ws.Range("J6").Value = 1000
$ws.Range("K6").Value = 22000
$ws.Range("L6").Value = 24000
$ws.Range("M6").Value = 23000
$ws.Range("P6").Value = 920

$ws.Range("D7").Value = 44363
$ws.Range("J7").Value = 240
$ws.Range("K7").Value = 28000
$ws.Range("L7").Value = 30000
$ws.Range("M7").Value = 29000
$ws.Range("P7").Value = 1160

$ws.Range("D8").Value = 44391
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 26000
$ws.Range("L8").Value = 28000
$ws.Range("M8").Value = 27000
$ws.Range("P8").Value = 1080

$ws.Range("D9").Value = 44364
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 28000
$ws.Range("L9").Value = 30000
$ws.Range("M9").Value = 29000
$ws.Range("P9").Value = 1160

$ws.Range("D10").Value = 44434
$ws.Range("J10").Value = 500
$ws.Range("K10").Value = 28000
$ws.Range("L10").Value = 30000
$ws.Range("M10").Value = 29000
$ws.Range("P10").Value = 1160

$ws.Range("D11").Value = 44441
$ws.Range("J11").Value = 700
$ws.Range("K11").Value = 28000
$ws.Range("L11").Value = 30000
$ws.Range("M11").Value = 29000
$ws.Range("P11").Value = 1160

$ws.Range("D12").Value = 44413
$ws.Range("J12").Value = 700
$ws.Range("K12").Value = 26000
$ws.Range("L12").Value = 28000
$ws.Range("M12").Value = 27000
$ws.Range("P12").Value = 1080

$ws.Range("D13").Value = 44448
$ws.Range("J13").Value = 400
$ws.Range("K13").Value = 28000
$ws.Range("L13").Value = 30000
$ws.Range("M13").Value = 29000
$ws.Range("P13").Value = 1160

$ws.Range("D14").Value = 44426
$ws.Range("J14").Value = 400
$ws.Range("K14").Value = 28000
$ws.Range("L14").Value = 30000
$ws.Range("M14").Value = 29000
$ws.Range("P14").Value = 1160

$ws.Range("D15").Value = 44489
$ws.Range("J15").Value = 400
$ws.Range("K15").Value = 18000
$ws.Range("L15").Value = 20000
$ws.Range("M15").Value = 19000
$ws.Range("P15").Value = 760

$ws.Range("D16").Value = 44461
$ws.Range("J16").Value = 500
$ws.Range("K16").Value = 23000
$ws.Range("L16").Value = 25000
$ws.Range("M16").Value = 24000
$ws.Range("P16").Value = 960

$ws.Range("D17").Value = 44455
$ws.Range("J17").Value = 800
$ws.Range("K17").Value = 28000
$ws.Range("L17").Value = 30000
$ws.Range("M17").Value = 29000
$ws.Range("P17").Value = 1160

$ws.Range("D18").Value = 44447
$ws.Range("J18").Value = 600
$ws.Range("K18").Value = 28000
$ws.Range("L18").Value = 30000
$ws.Range("M18").Value = 29000
$ws.Range("P18").Value = 1160

$ws.Range("D19").Value = 44483
$ws.Range("J19").Value = 300
$ws.Range("K19").Value = 18000
$ws.Range("L19").Value = 20000
$ws.Range("M19").Value = 19000
$ws.Range("P19").Value = 760

$ws.Range("D20").Value = 44435
$ws.Range("J20").Value = 900
$ws.Range("K20").Value = 28000
$ws.Range("L20").Value = 30000
$ws.Range("M20").Value = 29000
$ws.Range("P20").Value = 1160

$ws.Range("D21").Value = 44427
$ws.Range("J21").Value = 300
$ws.Range("K21").Value = 28000
$ws.Range("L21").Value = 30000
$ws.Range("M21").Value = 29000
$ws.Range("P21").Value = 1160

$ws.Range("D22").Value = 44468
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 23000
$ws.Range("L22").Value = 25000
$ws.Range("M22").Value = 24000
$ws.Range("P22").Value = 960

$ws.Range("D23").Value = 44420
$ws.Range("J23").Value = 700
$ws.Range("K23").Value = 27000
$ws.Range("L23").Value = 29000
$ws.Range("M23").Value = 28000
$ws.Range("P23").Value = 1120

$ws.Range("D24").Value = 44412
$ws.Range("J24").Value = 600
$ws.Range("K24").Value = 25000
$ws.Range("L24").Value = 27000
$ws.Range("M24").Value = 26000
$ws.Range("P24").Value = 1040

$ws.Range("D25").Value = 44377
$ws.Range("J25").Value = 500
$ws.Range("K25").Value = 26000
$ws.Range("L25").Value = 28000
$ws.Range("M25").Value = 27000
$ws.Range("P25").Value = 1080

$ws.Range("D26").Value = 44405
$ws.Range("J26").Value = 500
$ws.Range("K26").Value = 26000
$ws.Range("L26").Value = 28000
$ws.Range("M26").Value = 27000
$ws.Range("P26").Value = 1080

$ws.Range("D27").Value = 44350
$ws.Range("J27").Value = 700
$ws.Range("K27").Value = 28000
$ws.Range("L27").Value = 30000
$ws.Range("M27").Value = 29000
$ws.Range("P27").Value = 1160

$ws.Range("D28").Value = 44385
$ws.Range("J28").Value = 500
$ws.Range("K28").Value = 26000
$ws.Range("L28").Value = 28000
$ws.Range("M28").Value = 27000
$ws.Range("P28").Value = 1080

$ws.Range("D29").Value = 44398
$ws.Range("J29").Value = 500
$ws.Range("K29").Value = 26000
$ws.Range("L29").Value = 28000
$ws.Range("M29").Value = 27000
$ws.Range("P29").Value = 1080

$ws.Range("D30").Value = 44371
$ws.Range("J30").Value = 500
$ws.Range("K30").Value = 28000
$ws.Range("L30").Value = 30000
$ws.Range("M30").Value = 29000
$ws.Range("P30").Value = 1160

$ws.Range("D31").Value = 44454
$ws.Range("J31").Value = 1000
$ws.Range("K31").Value = 28000
$ws.Range("L31").Value = 30000
$ws.Range("M31").Value = 29000
$ws.Range("P31").Value = 1160

$ws.Range("D32").Value = 44490
$ws.Range("J32").Value = 500
$ws.Range("K32").Value = 16000
$ws.Range("L32").Value = 18000
$ws.Range("M32").Value = 17000
$ws.Range("P32").Value = 680

$ws.Range("D33").Value = 44399
$ws.Range("J33").Value = 400
$ws.Range("K33").Value = 26000
$ws.Range("L33").Value = 28000
$ws.Range("M33").Value = 27000
$ws.Range("P33").Value = 1080

$ws.Range("D34").Value = 44357
$ws.Range("J34").Value = 340
$ws.Range("K34").Value = 28000
$ws.Range("L34").Value = 30000
$ws.Range("M34").Value = 29000
$ws.Range("P34").Value = 1160

$ws.Range("D35").Value = 44476
$ws.Range("J35").Value = 500
$ws.Range("K35").Value = 23000
$ws.Range("L35").Value = 24000
$ws.Range("M35").Value = 23500
$ws.Range("P35").Value = 940

$ws.Range("D36").Value = 44482
$ws.Range("J36").Value = 500
$ws.Range("K36").Value = 18000
$ws.Range("L36").Value = 20000
$ws.Range("M36").Value = 19000
$ws.Range("P36").Value = 760

$ws.Range("D37").Value = 44356
$ws.Range("J37").Value = 300
$ws.Range("K37").Value = 26000
$ws.Range("L37").Value = 28000
$ws.Range("M37").Value = 27000
$ws.Range("P37").Value = 1080

$ws.Range("D38").Value = 44469
$ws.Range("J38").Value = 600
$ws.Range("K38").Value = 22000
$ws.Range("L38").Value = 24000
$ws.Range("M38").Value = 23000
$ws.Range("P38").Value = 920

$ws.Range("D40").Value = 44343
$ws.Range("J40").Value = 200
$ws.Range("K40").Value = 26000
$ws.Range("L40").Value = 28000
$ws.Range("M40").Value = 27000
$ws.Range("P40").Value = 1080

$ws.Range("D41").Value = 44370
$ws.Range("J41").Value = 400
$ws.Range("K41").Value = 27000
$ws.Range("L41").Value = 28000
$ws.Range("M41").Value = 27500
$ws.Range("P41").Value = 1100

$ws.Range("D43").Value = 44433
$ws.Range("J43").Value = 400
$ws.Range("K43").Value = 28000
$ws.Range("L43").Value = 30000
$ws.Range("M43").Value = 29000
$ws.Range("P43").Value = 1160
